# Generate Report for Handoff
#
# The localization status report is regenerated. The file
# "b9d59381-7ed0-45c7-9e7a-f365bb11980c" (row 6 on both the "zh-cn" and
# "de-de" sheets) has just been handed off again, so its
# "Latest Handoff Datetime" (column D) is refreshed with a new timestamp.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D6").Value = "2016-03-10 14:32:15"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D6").Value = "2016-03-10 14:32:19"
